$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 for the new "FAQ" test case, shifting the
# existing rows (old row 13 "Contact_Support" and below) down by one.
$ws.Rows("13:13").Insert()

# Match the formatting (borders/fill) used by the rest of the data rows.
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false | Out-Null

# Populate the new FAQ row (A13/B13/C13)
$ws.Range("A13").Value = "FAQ"
$ws.Range("C13").Value = "N"

# Contact_Support (now row 14) Runmode flips from Y to N
$ws.Range("C14").Value = "N"

# SendFeedback (now row 19) Runmode flips from N to Y
$ws.Range("C19").Value = "Y"

# Keep the active selection consistent with the authored workbook state
$ws.Range("C19").Select() | Out-Null
